$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1012503.44
$ws.Range("J17").Value = 1012503.44
$ws.Range("L17").Value = 3037510.32
$ws.Range("N17").Value = -3037846.32
$ws.Range("H18").Value = 476.55554
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H86").Value = 3098
$ws.Range("I86").Value = 3084.5
$ws.Range("J86").Value = 3107
$ws.Range("K86").Value = 3084.5
$ws.Range("L86").Value = 3107
$ws.Range("M86").Value = -1961.5
$ws.Range("N86").Value = -5353
$ws.Range("H89").Value = 3098
$ws.Range("I89").Value = 3084.5
$ws.Range("J89").Value = 3107
$ws.Range("K89").Value = 15422.5
$ws.Range("L89").Value = 15535
$ws.Range("M89").Value = -9806.5
$ws.Range("N89").Value = -26767
$ws.Range("H92").Value = 545.7143
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H96").Value = 2640.7368
$ws.Range("I96").Value = 2206.7273
$ws.Range("J96").Value = 3237.5
$ws.Range("K96").Value = 6620.1819
$ws.Range("L96").Value = 9712.5
$ws.Range("M96").Value = -5247.1819
$ws.Range("N96").Value = -12458.5
$ws.Range("H137").Value = 5200.385
$ws.Range("I137").Value = 5830.143
$ws.Range("J137").Value = 4465.6665
$ws.Range("K137").Value = 17490.429
$ws.Range("L137").Value = 13396.9995
$ws.Range("M137").Value = -14940.429
$ws.Range("N137").Value = -18496.9995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5420.1357
$ws.Range("I32").Value = 4542.673
$ws.Range("K32").Value = 4542.673
$ws.Range("M32").Value = -4255.673
$ws.Range("H69").Value = 180459
$ws.Range("J69").Value = 180459
$ws.Range("L69").Value = 180459
$ws.Range("N69").Value = -181957
$ws.Range("H72").Value = 180459
$ws.Range("J72").Value = 180459
$ws.Range("L72").Value = 541377
$ws.Range("N72").Value = -548865
$ws.Range("H132").Value = 33389008
$ws.Range("I132").Value = 12754.2
$ws.Range("J132").Value = 100141510
$ws.Range("K132").Value = 38262.60000000001
$ws.Range("L132").Value = 300424530
$ws.Range("M132").Value = -35732.60000000001
$ws.Range("N132").Value = -300429590

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 180000
$ws.Range("J70").Value = 180000
$ws.Range("L70").Value = 180000
$ws.Range("N70").Value = -180586
$ws.Range("H73").Value = 180000
$ws.Range("J73").Value = 180000
$ws.Range("L73").Value = 180000
$ws.Range("N73").Value = -182028
$ws.Range("H134").Value = 3239.65
$ws.Range("I134").Value = 3268.6924
$ws.Range("K134").Value = 9806.0772
$ws.Range("M134").Value = -7271.0772

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 500005000
$ws.Range("I35").Value = 500005000
$ws.Range("K35").Value = 500005000
$ws.Range("M35").Value = -500004706

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 105
$ws.Range("J2").Value = 105
$ws.Range("L2").Value = 630
$ws.Range("N2").Value = -856
$ws.Range("H39").Value = 763.7619
$ws.Range("I39").Value = 383.75
$ws.Range("J39").Value = 1979.8
$ws.Range("K39").Value = 1151.25
$ws.Range("L39").Value = 5939.4
$ws.Range("M39").Value = -857.25
$ws.Range("N39").Value = -6527.4
$ws.Range("H55").Value = 742
$ws.Range("J55").Value = 914.8570999999999
$ws.Range("L55").Value = 2744.5713
$ws.Range("N55").Value = -3098.5713
$ws.Range("H64").Value = 7875.316
$ws.Range("I64").Value = 2261
$ws.Range("J64").Value = 9880.429
$ws.Range("K64").Value = 6783
$ws.Range("L64").Value = 29641.287
$ws.Range("M64").Value = -6513
$ws.Range("N64").Value = -30181.287
$ws.Range("H67").Value = 7875.316
$ws.Range("I67").Value = 2261
$ws.Range("J67").Value = 9880.429
$ws.Range("K67").Value = 6783
$ws.Range("L67").Value = 29641.287
$ws.Range("M67").Value = -5847
$ws.Range("N67").Value = -31513.287
$ws.Range("H74").Value = 54000
$ws.Range("J74").Value = 54000
$ws.Range("L74").Value = 162000
$ws.Range("N74").Value = -164122
$ws.Range("H77").Value = 54000
$ws.Range("J77").Value = 54000
$ws.Range("L77").Value = 486000
$ws.Range("N77").Value = -496608
$ws.Range("H97").Value = 262.8889
$ws.Range("J97").Value = 260
$ws.Range("L97").Value = 780
$ws.Range("N97").Value = -1772
$ws.Range("H112").Value = 5000
$ws.Range("I112").Value = 5000
$ws.Range("K112").Value = 15000
$ws.Range("M112").Value = -13892
$ws.Range("H125").Value = 8500
$ws.Range("I125").Value = 7000
$ws.Range("J125").Value = 10000
$ws.Range("K125").Value = 21000
$ws.Range("L125").Value = 30000
$ws.Range("M125").Value = -16080
$ws.Range("N125").Value = -39840

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 31086.078
$ws.Range("J131").Value = 4663.6895
$ws.Range("L131").Value = 13991.0685
$ws.Range("N131").Value = -24071.0685
$ws.Range("H21").Value = 40000
$ws.Range("I21").Value = 40000
$ws.Range("K21").Value = 40000
$ws.Range("M21").Value = -39827
$ws.Range("H22").Value = 1350
$ws.Range("I22").Value = 1200
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -671
$ws.Range("N22").Value = -2558
$ws.Range("H24").Value = 56363.09
$ws.Range("I24").Value = 20000
$ws.Range("K24").Value = 20000
$ws.Range("M24").Value = -19827
$ws.Range("H30").Value = 40000
$ws.Range("I30").Value = 40000
$ws.Range("K30").Value = 40000
$ws.Range("M30").Value = -39895

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 23123.5
$ws.Range("I56").Value = 15999.667
$ws.Range("K56").Value = 15999.667
$ws.Range("M56").Value = -15308.667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10080.493
$ws.Range("I132").Value = 7816.814
$ws.Range("J132").Value = 13325.1
$ws.Range("K132").Value = 23450.442
$ws.Range("L132").Value = 39975.3
$ws.Range("M132").Value = -20920.442
$ws.Range("N132").Value = -45035.3
$ws.Range("H133").Value = 69447.5
$ws.Range("J133").Value = 69447.5
$ws.Range("L133").Value = 69447.5
$ws.Range("N133").Value = -74507.5
$ws.Range("H122").Value = 69223.87
$ws.Range("I122").Value = 85447.414
$ws.Range("J122").Value = 4329.6665
$ws.Range("K122").Value = 256342.242
$ws.Range("L122").Value = 13991.0685
$ws.Range("M122").Value = -253892.242
$ws.Range("N122").Value = -17888.9995
$ws.Range("H126").Value = 5382.4287
$ws.Range("I126").Value = 6021.1665
$ws.Range("K126").Value = 18063.4995
$ws.Range("M126").Value = -15593.4995
